# Update each answer cell in the 20x5 practice table in document order.
# Cell.Range.Text assignment is used (rather than Find/Replace) because it
# targets the specific cell by position, which matters here since several
# old/new values collide across different cells (e.g. "26x80=2080" appears
# both as an old value in one cell and a new value in another).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Cell(1,1): "48×38=1824" -> "96×71=6816"
$cell = $t.Cell(1, 1)
$cell.Range.Text = "96×71=6816"
# Cell(1,2): "54×19=1026" -> "97×93=9021"
$cell = $t.Cell(1, 2)
$cell.Range.Text = "97×93=9021"
# Cell(1,3): "45×99=4455" -> "98×48=4704"
$cell = $t.Cell(1, 3)
$cell.Range.Text = "98×48=4704"
# Cell(1,4): "12×57=684" -> "57×53=3021"
$cell = $t.Cell(1, 4)
$cell.Range.Text = "57×53=3021"
# Cell(1,5): "54×73=3942" -> "44×15=660"
$cell = $t.Cell(1, 5)
$cell.Range.Text = "44×15=660"
# Cell(5,1): "62×80=4960" -> "54×69=3726"
$cell = $t.Cell(5, 1)
$cell.Range.Text = "54×69=3726"
# Cell(5,2): "67×31=2077" -> "45×34=1530"
$cell = $t.Cell(5, 2)
$cell.Range.Text = "45×34=1530"
# Cell(5,3): "84×47=3948" -> "17×93=1581"
$cell = $t.Cell(5, 3)
$cell.Range.Text = "17×93=1581"
# Cell(5,4): "71×24=1704" -> "45×30=1350"
$cell = $t.Cell(5, 4)
$cell.Range.Text = "45×30=1350"
# Cell(5,5): "21×90=1890" -> "26×80=2080"
$cell = $t.Cell(5, 5)
$cell.Range.Text = "26×80=2080"
# Cell(10,1): "43×72=3096" -> "61×91=5551"
$cell = $t.Cell(10, 1)
$cell.Range.Text = "61×91=5551"
# Cell(10,2): "26×80=2080" -> "98×15=1470"
$cell = $t.Cell(10, 2)
$cell.Range.Text = "98×15=1470"
# Cell(10,3): "55×62=3410" -> "88×98=8624"
$cell = $t.Cell(10, 3)
$cell.Range.Text = "88×98=8624"
# Cell(10,4): "70×40=2800" -> "77×99=7623"
$cell = $t.Cell(10, 4)
$cell.Range.Text = "77×99=7623"
# Cell(10,5): "50×72=3600" -> "33×59=1947"
$cell = $t.Cell(10, 5)
$cell.Range.Text = "33×59=1947"
# Cell(15,1): "34×71=2414" -> "85×68=5780"
$cell = $t.Cell(15, 1)
$cell.Range.Text = "85×68=5780"
# Cell(15,2): "16×71=1136" -> "80×16=1280"
$cell = $t.Cell(15, 2)
$cell.Range.Text = "80×16=1280"
# Cell(15,3): "67×95=6365" -> "57×24=1368"
$cell = $t.Cell(15, 3)
$cell.Range.Text = "57×24=1368"
# Cell(15,4): "50×71=3550" -> "42×81=3402"
$cell = $t.Cell(15, 4)
$cell.Range.Text = "42×81=3402"
# Cell(15,5): "16×92=1472" -> "69×80=5520"
$cell = $t.Cell(15, 5)
$cell.Range.Text = "69×80=5520"
# Cell(20,1): "42×52=2184" -> "78×12=936"
$cell = $t.Cell(20, 1)
$cell.Range.Text = "78×12=936"
# Cell(20,2): "33×27=891" -> "25×21=525"
$cell = $t.Cell(20, 2)
$cell.Range.Text = "25×21=525"
# Cell(20,3): "19×26=494" -> "70×55=3850"
$cell = $t.Cell(20, 3)
$cell.Range.Text = "70×55=3850"
# Cell(20,4): "21×88=1848" -> "78×44=3432"
$cell = $t.Cell(20, 4)
$cell.Range.Text = "78×44=3432"
# Cell(20,5): "41×36=1476" -> "92×74=6808"
$cell = $t.Cell(20, 5)
$cell.Range.Text = "92×74=6808"
